$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark (it will be re-added later at its new location).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. "Researched..." bullet: split wording into two sentences/runs.
$d.Content.Find.Execute(
    "Researched the various possible applications that we can develop",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Researched the various possible approaches to develop the application", 2) | Out-Null

# 3. Merge "Started creating API for " + "add account" + " and " runs.
$d.Content.Find.Execute(
    "Started creating API for add account and ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Started creating API for add account and ", 2) | Out-Null

# 4. Merge "Started creating API for " + "transfer between accounts" + " " runs.
$d.Content.Find.Execute(
    "Started creating API for transfer between accounts ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Started creating API for transfer between accounts ", 2) | Out-Null

# 5. Merge "View and search Transactions " + "–" + " for" + " " + "credits/debits/checks/fees".
$d.Content.Find.Execute(
    "View and search Transactions – for credits/debits/checks/fees",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "View and search Transactions – for credits/debits/checks/fees", 2) | Out-Null

# 6. Merge "C" + "reat" + "ed" + " API for login and " runs -> "Created API for login and ".
$d.Content.Find.Execute(
    "Created API for login and ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created API for login and ", 2) | Out-Null

# 7. Merge "C" + "reat" + "ed" + " API for add account and " runs -> "Created API for add account and ".
$d.Content.Find.Execute(
    "Created API for add account and ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created API for add account and ", 2) | Out-Null

# 8. Merge "C" + "reat" + "ed" + " API for transfer between accounts " runs, plus merge the
#    "- one time or recurring" + " " runs right after it.
$d.Content.Find.Execute(
    "Created API for transfer between accounts ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created API for transfer between accounts ", 2) | Out-Null

$d.Content.Find.Execute(
    "- one time or recurring ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "- one time or recurring ", 2) | Out-Null

# 9. Merge "C" + "reat" + "ed" + " API for " runs -> "Created API for ".
$d.Content.Find.Execute(
    "Created API for ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created API for ", 2) | Out-Null

# 10. Re-add the _GoBack bookmark, now collapsed right after the "Researched..." bullet text.
$findRange = $d.Content
$findRange.Find.Execute(
    "Researched the various possible approaches to develop the application",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $findRange) | Out-Null
